# #5: property boat&car done
# Rebuild the "汽車" (car) sheet (3rd worksheet) so row 1 becomes a proper
# header row and row 2 gains the property_category / category / date /
# legislator_name / legislator_id / source_file / index columns already
# present on the 土地 (land) and 建物 (building) sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: turn the old (duplicated-data) header row into real headers ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the header styling (bold + border, like B1 already had) across the
# newly-used header cells H1:N1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null

# --- Row 2: keep existing A2:G2 data, append the new tracked columns ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# Leading apostrophe keeps this a text cell instead of Excel auto-parsing it
# as a date serial number.
$ws.Range("J2").Value = "'2012-04-26"
$ws.Range("K2").Value = "楊瓊瓔"
$ws.Range("L2").Value = 854
$ws.Range("M2").Value = "tmp8a701"
$ws.Range("N2").Value = 44
